# Apply the timesheet edits described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the description text in B1 (truncated text -> full text)
$ws.Range("B1").Value = "To engage Agile Co-Development and ICT Professional Services via 19024 bulk tender (PR230941)"

# 2. Rows 12 & 13: the "Annual Leave" (G) mark moves to "At Work" (C)
#    Values such as "1.0" look numeric, so a leading apostrophe is used to
#    keep them stored as literal text (matching the original inlineStr cells).
#    A lone apostrophe produces an empty *text* cell (as opposed to assigning
#    "" which would null the cell out and change its stored type).
$ws.Range("C12").Value = "'1.0"
$ws.Range("G12").Value = "'"

$ws.Range("C13").Value = "'1.0"
$ws.Range("G13").Value = "'"

# 3. Rows 32, 33, 34, 37, 38, 41: the "At Work" (C) mark moves to "Annual Leave" (G)
$ws.Range("C32").Value = "'"
$ws.Range("G32").Value = "'1.0"

$ws.Range("C33").Value = "'"
$ws.Range("G33").Value = "'1.0"

$ws.Range("C34").Value = "'"
$ws.Range("G34").Value = "'1.0"

$ws.Range("C37").Value = "'"
$ws.Range("G37").Value = "'1.0"

$ws.Range("C38").Value = "'"
$ws.Range("G38").Value = "'1.0"

$ws.Range("C41").Value = "'"
$ws.Range("G41").Value = "'1.0"

# 4. Update totals row 44: At Work total 18.0 -> 14.0, Annual Leave total 2.0 -> 6.0
$ws.Range("C44").Value = "'14.0"
$ws.Range("G44").Value = "'6.0"
